$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"235.2675016666667"
$ws.Range("H2").Value = [double]"705.802505"
$ws.Range("I2").Value = [double]"0.5738994362335403"
$ws.Range("J2").Value = [double]"0.5738994362335402"
$ws.Range("M2").Value = [double]"6.066605666666667"
$ws.Range("N2").Value = [double]"18.199817"
$ws.Range("O2").Value = [double]"0.8497846287916651"
$ws.Range("P2").Value = [double]"0.8497846287916652"
$ws.Range("Q2").Value = [double]"1427.275158793509"
$ws.Range("R2").Value = [double]"12845.47642914158"
$ws.Range("S2").Value = [double]"0.4876909193834649"
$ws.Range("T2").Value = [double]"0.4876909193834649"
$ws.Range("G3").Value = [double]"235.2675016666667"
$ws.Range("H3").Value = [double]"705.802505"
$ws.Range("I3").Value = [double]"0.5738994362335403"
$ws.Range("J3").Value = [double]"0.5738994362335402"
$ws.Range("O3").Value = [double]"0.1196497582104962"
$ws.Range("P3").Value = [double]"0.1196497582104962"
$ws.Range("Q3").Value = [double]"200.9604808836311"
$ws.Range("R3").Value = [double]"1808.64432795268"
$ws.Range("S3").Value = [double]"0.06866692878248316"
$ws.Range("T3").Value = [double]"0.06866692878248316"
$ws.Range("G4").Value = [double]"235.2675016666667"
$ws.Range("H4").Value = [double]"705.802505"
$ws.Range("I4").Value = [double]"0.5738994362335403"
$ws.Range("J4").Value = [double]"0.5738994362335402"
$ws.Range("M4").Value = [double]"0.1824346666666667"
$ws.Range("N4").Value = [double]"0.547304"
$ws.Range("O4").Value = [double]"0.02555468148257719"
$ws.Range("P4").Value = [double]"0.02555468148257719"
$ws.Range("Q4").Value = [double]"42.92094824405778"
$ws.Range("R4").Value = [double]"386.28853419652"
$ws.Range("S4").Value = [double]"0.01466581729597874"
$ws.Range("T4").Value = [double]"0.01466581729597874"
$ws.Range("G5").Value = [double]"235.2675016666667"
$ws.Range("H5").Value = [double]"705.802505"
$ws.Range("I5").Value = [double]"0.5738994362335403"
$ws.Range("J5").Value = [double]"0.5738994362335402"
$ws.Range("M5").Value = [double]"0.035773"
$ws.Range("N5").Value = [double]"0.107319"
$ws.Range("O5").Value = [double]"0.005010931515261538"
$ws.Range("P5").Value = [double]"0.005010931515261539"
$ws.Range("Q5").Value = [double]"8.416224337121667"
$ws.Range("R5").Value = [double]"75.746019034095"
$ws.Range("S5").Value = [double]"0.002875770771613477"
$ws.Range("T5").Value = [double]"0.002875770771613477"
$ws.Range("I6").Value = [double]"0.3286113026040369"
$ws.Range("J6").Value = [double]"0.3286113026040369"
$ws.Range("M6").Value = [double]"6.066605666666667"
$ws.Range("N6").Value = [double]"18.199817"
$ws.Range("O6").Value = [double]"0.8497846287916651"
$ws.Range("P6").Value = [double]"0.8497846287916652"
$ws.Range("Q6").Value = [double]"817.2490152345404"
$ws.Range("R6").Value = [double]"7355.241137110864"
$ws.Range("S6").Value = [double]"0.2792488338001171"
$ws.Range("T6").Value = [double]"0.279248833800117"
$ws.Range("I7").Value = [double]"0.3286113026040369"
$ws.Range("J7").Value = [double]"0.3286113026040369"
$ws.Range("O7").Value = [double]"0.1196497582104962"
$ws.Range("P7").Value = [double]"0.1196497582104962"
$ws.Range("S7").Value = [double]"0.03931826290180921"
$ws.Range("T7").Value = [double]"0.03931826290180921"
$ws.Range("I8").Value = [double]"0.3286113026040369"
$ws.Range("J8").Value = [double]"0.3286113026040369"
$ws.Range("M8").Value = [double]"0.1824346666666667"
$ws.Range("N8").Value = [double]"0.547304"
$ws.Range("O8").Value = [double]"0.02555468148257719"
$ws.Range("P8").Value = [double]"0.02555468148257719"
$ws.Range("Q8").Value = [double]"24.57627211492978"
$ws.Range("R8").Value = [double]"221.186449034368"
$ws.Range("S8").Value = [double]"0.008397557169620951"
$ws.Range("T8").Value = [double]"0.008397557169620951"
$ws.Range("I9").Value = [double]"0.3286113026040369"
$ws.Range("J9").Value = [double]"0.3286113026040369"
$ws.Range("M9").Value = [double]"0.035773"
$ws.Range("N9").Value = [double]"0.107319"
$ws.Range("O9").Value = [double]"0.005010931515261538"
$ws.Range("P9").Value = [double]"0.005010931515261539"
$ws.Range("Q9").Value = [double]"4.819078514138666"
$ws.Range("R9").Value = [double]"43.371706627248"
$ws.Range("S9").Value = [double]"0.001646648732489715"
$ws.Range("T9").Value = [double]"0.001646648732489715"
$ws.Range("G10").Value = [double]"0.325805"
$ws.Range("H10").Value = [double]"0.977415"
$ws.Range("I10").Value = [double]"0.0007947519504286909"
$ws.Range("J10").Value = [double]"0.0007947519504286907"
$ws.Range("M10").Value = [double]"6.066605666666667"
$ws.Range("N10").Value = [double]"18.199817"
$ws.Range("O10").Value = [double]"0.8497846287916651"
$ws.Range("P10").Value = [double]"0.8497846287916652"
$ws.Range("Q10").Value = [double]"1.976530459228333"
$ws.Range("R10").Value = [double]"17.788774133055"
$ws.Range("S10").Value = [double]"0.000675367991176497"
$ws.Range("T10").Value = [double]"0.0006753679911764968"
$ws.Range("G11").Value = [double]"0.325805"
$ws.Range("H11").Value = [double]"0.977415"
$ws.Range("I11").Value = [double]"0.0007947519504286909"
$ws.Range("J11").Value = [double]"0.0007947519504286907"
$ws.Range("O11").Value = [double]"0.1196497582104962"
$ws.Range("P11").Value = [double]"0.1196497582104962"
$ws.Range("Q11").Value = [double]"0.2782956804933333"
$ws.Range("R11").Value = [double]"2.50466112444"
$ws.Range("S11").Value = [double]"9.509187870611311E-05"
$ws.Range("T11").Value = [double]"9.509187870611311E-05"
$ws.Range("G12").Value = [double]"0.325805"
$ws.Range("H12").Value = [double]"0.977415"
$ws.Range("I12").Value = [double]"0.0007947519504286909"
$ws.Range("J12").Value = [double]"0.0007947519504286907"
$ws.Range("M12").Value = [double]"0.1824346666666667"
$ws.Range("N12").Value = [double]"0.547304"
$ws.Range("O12").Value = [double]"0.02555468148257719"
$ws.Range("P12").Value = [double]"0.02555468148257719"
$ws.Range("Q12").Value = [double]"0.05943812657333333"
$ws.Range("R12").Value = [double]"0.53494313916"
$ws.Range("S12").Value = [double]"2.030963295086217E-05"
$ws.Range("T12").Value = [double]"2.030963295086217E-05"
$ws.Range("G13").Value = [double]"0.325805"
$ws.Range("H13").Value = [double]"0.977415"
$ws.Range("I13").Value = [double]"0.0007947519504286909"
$ws.Range("J13").Value = [double]"0.0007947519504286907"
$ws.Range("M13").Value = [double]"0.035773"
$ws.Range("N13").Value = [double]"0.107319"
$ws.Range("O13").Value = [double]"0.005010931515261538"
$ws.Range("P13").Value = [double]"0.005010931515261539"
$ws.Range("Q13").Value = [double]"0.011655022265"
$ws.Range("R13").Value = [double]"0.104895200385"
$ws.Range("S13").Value = [double]"3.982447595218703E-06"
$ws.Range("T13").Value = [double]"3.982447595218702E-06"
$ws.Range("G14").Value = [double]"39.46134166666666"
$ws.Range("H14").Value = [double]"118.384025"
$ws.Range("I14").Value = [double]"0.09625996610278018"
$ws.Range("J14").Value = [double]"0.09625996610278018"
$ws.Range("M14").Value = [double]"6.066605666666667"
$ws.Range("N14").Value = [double]"18.199817"
$ws.Range("O14").Value = [double]"0.8497846287916651"
$ws.Range("P14").Value = [double]"0.8497846287916652"
$ws.Range("Q14").Value = [double]"239.3963989692694"
$ws.Range("R14").Value = [double]"2154.567590723425"
$ws.Range("S14").Value = [double]"0.08180023956214931"
$ws.Range("T14").Value = [double]"0.08180023956214932"
$ws.Range("G15").Value = [double]"39.46134166666666"
$ws.Range("H15").Value = [double]"118.384025"
$ws.Range("I15").Value = [double]"0.09625996610278018"
$ws.Range("J15").Value = [double]"0.09625996610278018"
$ws.Range("O15").Value = [double]"0.1196497582104962"
$ws.Range("P15").Value = [double]"0.1196497582104962"
$ws.Range("Q15").Value = [double]"33.7070362097111"
$ws.Range("R15").Value = [double]"303.3633258873999"
$ws.Range("S15").Value = [double]"0.0115174816695482"
$ws.Range("T15").Value = [double]"0.01151748166954821"
$ws.Range("G16").Value = [double]"39.46134166666666"
$ws.Range("H16").Value = [double]"118.384025"
$ws.Range("I16").Value = [double]"0.09625996610278018"
$ws.Range("J16").Value = [double]"0.09625996610278018"
$ws.Range("M16").Value = [double]"0.1824346666666667"
$ws.Range("N16").Value = [double]"0.547304"
$ws.Range("O16").Value = [double]"0.02555468148257719"
$ws.Range("P16").Value = [double]"0.02555468148257719"
$ws.Range("Q16").Value = [double]"7.199116713177776"
$ws.Range("R16").Value = [double]"64.7920504186"
$ws.Range("S16").Value = [double]"0.002459892773280224"
$ws.Range("T16").Value = [double]"0.002459892773280225"
$ws.Range("G17").Value = [double]"39.46134166666666"
$ws.Range("H17").Value = [double]"118.384025"
$ws.Range("I17").Value = [double]"0.09625996610278018"
$ws.Range("J17").Value = [double]"0.09625996610278018"
$ws.Range("M17").Value = [double]"0.035773"
$ws.Range("N17").Value = [double]"0.107319"
$ws.Range("O17").Value = [double]"0.005010931515261538"
$ws.Range("P17").Value = [double]"0.005010931515261539"
$ws.Range("Q17").Value = [double]"1.411650575441666"
$ws.Range("R17").Value = [double]"12.704855178975"
$ws.Range("S17").Value = [double]"0.0004823520978024286"
$ws.Range("T17").Value = [double]"0.0004823520978024286"
$ws.Range("G18").Value = [double]"0.178139"
$ws.Range("H18").Value = [double]"0.534417"
$ws.Range("I18").Value = [double]"0.0004345431092138444"
$ws.Range("J18").Value = [double]"0.0004345431092138443"
$ws.Range("M18").Value = [double]"6.066605666666667"
$ws.Range("N18").Value = [double]"18.199817"
$ws.Range("O18").Value = [double]"0.8497846287916651"
$ws.Range("P18").Value = [double]"0.8497846287916652"
$ws.Range("Q18").Value = [double]"1.080699066854333"
$ws.Range("R18").Value = [double]"9.726291601689001"
$ws.Range("S18").Value = [double]"0.0003692680547572628"
$ws.Range("T18").Value = [double]"0.0003692680547572627"
$ws.Range("G19").Value = [double]"0.178139"
$ws.Range("H19").Value = [double]"0.534417"
$ws.Range("I19").Value = [double]"0.0004345431092138444"
$ws.Range("J19").Value = [double]"0.0004345431092138443"
$ws.Range("O19").Value = [double]"0.1196497582104962"
$ws.Range("P19").Value = [double]"0.1196497582104962"
$ws.Range("Q19").Value = [double]"0.1521625335013333"
$ws.Range("R19").Value = [double]"1.369462801512"
$ws.Range("S19").Value = [double]"5.199297794947371E-05"
$ws.Range("T19").Value = [double]"5.199297794947371E-05"
$ws.Range("G20").Value = [double]"0.178139"
$ws.Range("H20").Value = [double]"0.534417"
$ws.Range("I20").Value = [double]"0.0004345431092138444"
$ws.Range("J20").Value = [double]"0.0004345431092138443"
$ws.Range("M20").Value = [double]"0.1824346666666667"
$ws.Range("N20").Value = [double]"0.547304"
$ws.Range("O20").Value = [double]"0.02555468148257719"
$ws.Range("P20").Value = [double]"0.02555468148257719"
$ws.Range("Q20").Value = [double]"0.03249872908533333"
$ws.Range("R20").Value = [double]"0.292488561768"
$ws.Range("S20").Value = [double]"1.110461074640855E-05"
$ws.Range("T20").Value = [double]"1.110461074640854E-05"
$ws.Range("G21").Value = [double]"0.178139"
$ws.Range("H21").Value = [double]"0.534417"
$ws.Range("I21").Value = [double]"0.0004345431092138444"
$ws.Range("J21").Value = [double]"0.0004345431092138443"
$ws.Range("M21").Value = [double]"0.035773"
$ws.Range("N21").Value = [double]"0.107319"
$ws.Range("O21").Value = [double]"0.005010931515261538"
$ws.Range("P21").Value = [double]"0.005010931515261539"
$ws.Range("Q21").Value = [double]"0.006372566447000001"
$ws.Range("R21").Value = [double]"0.057353098023"
$ws.Range("S21").Value = [double]"2.177465760699389E-06"
$ws.Range("T21").Value = [double]"2.177465760699389E-06"
